# Generate Report for Handoff
# Adds a new tracked file "ebf26689-adde-4d6f-9f8a-595b0502a417.md" (status:
# "Ready for handoff") as row 9 on the Overview / zh-cn / de-de sheets, and
# grows each sheet's table to include it.

$wb = $excel.ActiveWorkbook

$fileName   = "ebf26689-adde-4d6f-9f8a-595b0502a417.md"
$pathName   = "e2e\" + $fileName
$ext        = ".md"
$status     = "Ready for handoff"
$hoDate     = "2016-08-31 12:01:57"
$zhHoDate   = "2016-08-31 12:01:47"
$zhXlf      = "ebf26689-adde-4d6f-9f8a-595b0502a417.ea595fcb5cc6c637b1a03b8612b6cf49fd88ef84.zh-cn.xlf"
$deXlf      = "ebf26689-adde-4d6f-9f8a-595b0502a417.ea595fcb5cc6c637b1a03b8612b6cf49fd88ef84.de-de.xlf"
$epoch      = "0001-01-01 00:00:00"

$baseUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea595fcb5cc6c637b1a03b8612b6cf49fd88ef84/e2e/" + $fileName

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value2 = $fileName
$wsOverview.Range("B9").Value2 = $pathName
$wsOverview.Range("C9").Value2 = $ext
$wsOverview.Range("D9").Value2 = ""
$wsOverview.Range("E9").Value2 = $status
$wsOverview.Range("F9").Value2 = $status
$wsOverview.Range("G9").Value2 = $hoDate
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $baseUrl, "", "", $pathName) | Out-Null
$wsOverview.Range("B9").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A9").Value2 = $fileName
$wsZh.Range("B9").Value2 = $ext
$wsZh.Range("C9").Value2 = $status
$wsZh.Range("D9").Value2 = "e2e"
$wsZh.Range("E9").Value2 = "ht"
$wsZh.Range("F9").Value2 = "'False"
$wsZh.Range("G9").Value2 = $zhXlf
$wsZh.Range("H9").Value2 = $zhHoDate
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I9").Value2 = ""
$wsZh.Range("J9").Value2 = ""
$wsZh.Range("K9").Value2 = $epoch
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L9").Value2 = ""
$wsZh.Range("M9").Value2 = "'True"
$wsZh.Range("N9").Value2 = ""
$wsZh.Range("O9").Value2 = "'False"
$wsZh.Range("P9").Value2 = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $baseUrl, "", "", $fileName) | Out-Null
$wsZh.Range("A9").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A9").Value2 = $fileName
$wsDe.Range("B9").Value2 = $ext
$wsDe.Range("C9").Value2 = $status
$wsDe.Range("D9").Value2 = "e2e"
$wsDe.Range("E9").Value2 = "ht"
$wsDe.Range("F9").Value2 = "'False"
$wsDe.Range("G9").Value2 = $deXlf
$wsDe.Range("H9").Value2 = $hoDate
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I9").Value2 = ""
$wsDe.Range("J9").Value2 = ""
$wsDe.Range("K9").Value2 = $epoch
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L9").Value2 = ""
$wsDe.Range("M9").Value2 = "'True"
$wsDe.Range("N9").Value2 = ""
$wsDe.Range("O9").Value2 = "'False"
$wsDe.Range("P9").Value2 = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $baseUrl, "", "", $fileName) | Out-Null
$wsDe.Range("A9").Style = "Hyperlink"

"Handoff row added to Overview, zh-cn, de-de"
